# "Update dan tambah draft article"
# Adds four new keyword rows (31-34) to Sheet1, mirroring the pattern of the
# existing rows: column A holds the keyword phrase, column B counts the
# number of words via a LEN/SUBSTITUTE formula, and column F is a manual
# "priority" flag (left blank on the last new row, same as the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New keyword rows -------------------------------------------------
$ws.Range("A31").Value = "contoh surat lamaran kerja bidan di rumah sakit swasta"
$ws.Range("A32").Value = "contoh surat lamaran kerja bidan di klinik bersalin"
$ws.Range("A33").Value = "contoh surat lamaran kerja di bidang kesehatan dalam bahasa inggris"
$ws.Range("A34").Value = "contoh surat lamaran kerja bidan untuk dinas kesehatan"

$ws.Range("B31").Formula = '=LEN(A31)-LEN(SUBSTITUTE(A31," ",""))+1'
$ws.Range("B32").Formula = '=LEN(A32)-LEN(SUBSTITUTE(A32," ",""))+1'
$ws.Range("B33").Formula = '=LEN(A33)-LEN(SUBSTITUTE(A33," ",""))+1'
$ws.Range("B34").Formula = '=LEN(A34)-LEN(SUBSTITUTE(A34," ",""))+1'

$ws.Range("F31").Value = 1
$ws.Range("F32").Value = 1
$ws.Range("F33").Value = 1

# Match the existing look of the table (column A / B formatting) by
# copying the formatting from the row directly above the new block.
$ws.Range("A30").Copy() | Out-Null
$ws.Range("A31:A34").PasteSpecial(-4122) | Out-Null

# --- Column A got a little wider to fit the new, longer keyword -------
$ws.Columns.Item(1).ColumnWidth = 62.25

# --- Scroll position / selection the file was left at -----------------
$ws.Range("J26").Select() | Out-Null
